# Swap the order of "System" and the email address in the
# "Recorded By" column (column G) wherever it appears as
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$searchRange = $ws.Columns.Item(7)  # Column G - "Recorded By"

$firstFoundAddress = $null
$cell = $searchRange.Find($target, [Type]::Missing, [Type]::Missing, [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

while ($cell -ne $null) {
    if ($firstFoundAddress -eq $null) {
        $firstFoundAddress = $cell.Address()
    } elseif ($cell.Address() -eq $firstFoundAddress) {
        break
    }

    $cell.Value2 = $replacement

    $cell = $searchRange.FindNext($cell)
}
